# Apply updated cryptocurrency market data (prices and 1h volume change %)
# Values are written as text (matching the original inlineStr cell type) by
# temporarily forcing a Text number format so Excel does not auto-convert
# numeric-looking strings (e.g. "1.00", "0.0314") into real numbers, then
# clearing the format again so no stray style is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

Set-TextValue "D2" "60.186.53"
Set-TextValue "E2" "  +2.69%  "
Set-TextValue "D3" "3.204.71"
Set-TextValue "E3" "  +1.28%  "
Set-TextValue "E4" "  +0.00%  "
Set-TextValue "D5" "538.99"
Set-TextValue "E5" "  +1.78%  "
Set-TextValue "D6" "145.85"
Set-TextValue "E6" "  +4.33%  "
Set-TextValue "D7" "1.00"
Set-TextValue "E7" "  -0.09%  "
Set-TextValue "D8" "0.529"
Set-TextValue "E8" "  -2.04%  "
Set-TextValue "E9" "  +0.68%  "
Set-TextValue "E10" "  +1.14%  "
Set-TextValue "D11" "0.434"
Set-TextValue "E11" "  -0.68%  "
Set-TextValue "D12" "3.758.41"
Set-TextValue "E12" "  +1.30%  "
Set-TextValue "E13" "  -1.96%  "
Set-TextValue "D14" "25.82"
Set-TextValue "E14" "  +0.39%  "
Set-TextValue "E15" "  +0.91%  "
Set-TextValue "D16" "60.205.36"
Set-TextValue "E16" "  +2.63%  "
Set-TextValue "D17" "3.199.65"
Set-TextValue "E17" "  +2.94%  "
Set-TextValue "D18" "6.31"
Set-TextValue "E18" "  +0.88%  "
Set-TextValue "D19" "13.26"
Set-TextValue "E19" "  +2.10%  "
Set-TextValue "D20" "8.22"
Set-TextValue "E20" "  +1.36%  "
Set-TextValue "D21" "370.94"
Set-TextValue "E21" "  -1.31%  "
Set-TextValue "E22" "  -0.02%  "
Set-TextValue "E23" "  -1.34%  "
Set-TextValue "D24" "69.62"
Set-TextValue "E24" "  -0.12%  "
Set-TextValue "E25" "  +1.84%  "
Set-TextValue "D26" "8.63"
Set-TextValue "E26" "  +4.45%  "
Set-TextValue "D27" "0.999"
Set-TextValue "E27" "  -0.31%  "
Set-TextValue "D28" "0.0₃0878"
Set-TextValue "E28" "  +1.44%  "
Set-TextValue "D29" "22.45"
Set-TextValue "E29" "  +0.41%  "
Set-TextValue "E30" "  +0.82%  "
Set-TextValue "E31" "  +1.59%  "
Set-TextValue "D32" "5.29"
Set-TextValue "E32" "  +2.28%  "
Set-TextValue "D33" "6.58"
Set-TextValue "E33" "  +4.28%  "
Set-TextValue "E34" "  +3.08%  "
Set-TextValue "D35" "159.03"
Set-TextValue "E35" "  +1.54%  "
Set-TextValue "E36" "  +2.87%  "
Set-TextValue "E37" "  +5.55%  "
Set-TextValue "D38" "2.795.34"
Set-TextValue "E38" "  +4.37%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D39" "0.0314"
Set-TextValue "E39" "  +8.89%  "
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D40" "0.0710"
Set-TextValue "E40" "  +2.05%  "
Set-TextValue "D41" "1.69"
Set-TextValue "E41" "  +0.48%  "
Set-TextValue "D42" "4.21"
Set-TextValue "E42" "  -1.46%  "
Set-TextValue "D43" "39.96"
Set-TextValue "E43" "  +2.09%  "
Set-TextValue "D44" "0.720"
Set-TextValue "E44" "  -0.48%  "
Set-TextValue "E45" "  +1.09%  "
Set-TextValue "D46" "3.244.73"
Set-TextValue "E46" "  +1.18%  "
Set-TextValue "E47" "  +0.39%  "
Set-TextValue "D48" "6.17"
Set-TextValue "D49" "20.71"
Set-TextValue "E49" "  +3.02%  "
Set-TextValue "D50" "0.797"
Set-TextValue "E50" "  +6.12%  "
Set-TextValue "E51" "  +0.05%  "
